# Cleaned dataset: remove unwanted columns for the analysis (flag via new
# "Use" column) and fill in the "Correct type" column for rows that were
# previously left blank (null -> should have been a concrete type / zero).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the previously-blank "Correct type" column (D) ------------------
# (written first so the shared-string table picks up "float" before the new
# "Use"/"NO" strings, matching the original author's edit order)
$dFillOrder = @(35,36,38,39,40,44,45,46,52,53,54,55,56,61,63,64,65)
$dFills = @{
    35 = "int"
    36 = "int"
    38 = "int"
    39 = "int"
    40 = "int"
    44 = "int"
    45 = "int"
    46 = "int"
    52 = "int"
    53 = "int"
    54 = "int"
    55 = "int"
    56 = "int"
    61 = "int"
    63 = "object"
    64 = "float"
    65 = "datetime"
}
foreach ($r in $dFillOrder) {
    $ws.Cells.Item($r, 4).Value = $dFills[$r]
}

# --- New column E: "Use" flag -------------------------------------------------
$ws.Cells.Item(1, 5).Value = "Use"

$eNoRows = @(2,3,7,14,21,23,24,25,26,27,28,29,30,31,32,33,34,35,39,40,41,42,44,45,46,51,61,62,63,65,66,67,69,70,71,72,73,74,75,76)
foreach ($r in $eNoRows) {
    $ws.Cells.Item($r, 5).Value = "NO"
}

# --- Apply an AutoFilter over the full used range -----------------------------
$rng = $ws.Range("A1:E81")
$rng.AutoFilter() | Out-Null

# --- Update the view: clear the frozen/top-left scroll cell, move selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("K55").Select() | Out-Null
